# Updated cryptos list with GitHub Actions
#
# The "Price" column (D) is free-form text scraped from the site (it mixes
# thousands-separated values like "63.070.26" with plain decimals like
# "415.42"), so every D value here is written with a leading apostrophe to
# force Excel to keep it as text instead of silently parsing it as a number.
# The "Volume(1h)" column (E) already carries padding spaces and a trailing
# "%" so it is never auto-converted and can be written as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new D (price) value, exactly as it should read once in the cell.
$dValues = @{
    2  = "62.915.11"
    3  = "3.471.98"
    5  = "415.42"
    6  = "131.04"
    7  = "0.630"
    9  = "0.732"
    10 = "0.152"
    11 = "42.81"
    12 = "9.77"
    13 = "0.0000224"
    14 = "4.041.37"
    16 = "20.56"
    17 = "3.477.65"
    18 = "12.69"
    20 = "62.897.13"
    21 = "468.21"
    22 = "90.75"
    23 = "3.30"
    24 = "13.23"
    25 = "10.62"
    26 = "3.34"
    27 = "33.51"
    29 = "7.57"
    30 = "12.15"
    32 = "0.168"
    34 = "41.05"
    35 = "0.999"
    36 = "58.31"
    37 = "0.0491"
    39 = "3.07"
    40 = "2.78"
    41 = "0.135"
    42 = "148.20"
    45 = "3.34"
    46 = "2.06"
    47 = "0.0₃0575"
    48 = "2.39"
    49 = "16.40"
    51 = "0.144"
}

# Row -> new E (1h volume/percentage) value, with the surrounding spaces
# preserved exactly as in the source data.
$eValues = @{
    2  = "  +1.69%  "
    3  = "  +1.92%  "
    4  = "  +0.33%  "
    5  = "  +1.52%  "
    6  = "  +1.86%  "
    7  = "  -1.25%  "
    9  = "  -0.50%  "
    10 = "  +6.96%  "
    11 = "  -1.89%  "
    12 = "  +4.48%  "
    13 = "  -0.14%  "
    14 = "  +2.40%  "
    15 = "  -0.23%  "
    16 = "  -3.80%  "
    17 = "  +1.29%  "
    18 = "  +1.10%  "
    19 = "  +0.23%  "
    20 = "  +1.67%  "
    21 = "  +4.00%  "
    22 = "  -1.28%  "
    23 = "  +2.97%  "
    24 = "  +0.09%  "
    25 = "  +14.00%  "
    26 = "  +1.00%  "
    27 = "  +0.78%  "
    28 = "  -0.01%  "
    29 = "  -0.90%  "
    30 = "  +0.91%  "
    31 = "  -1.38%  "
    32 = "  -1.13%  "
    33 = "  -1.48%  "
    34 = "  -3.45%  "
    35 = "  -0.04%  "
    36 = "  +8.07%  "
    37 = "  -2.73%  "
    38 = "  +0.26%  "
    39 = "  +3.66%  "
    40 = "  +7.03%  "
    41 = "  -0.87%  "
    42 = "  +2.77%  "
    43 = "  +2.18%  "
    44 = "  +0.84%  "
    45 = "  -1.73%  "
    46 = "  +2.74%  "
    47 = "  +31.62%  "
    48 = "  +11.07%  "
    49 = "  -1.68%  "
    50 = "  -0.91%  "
    51 = "  -4.89%  "
}

foreach ($row in $dValues.Keys) {
    $text = $dValues[$row]
    $cell = $ws.Range("D$row")
    # Force text entry (like typing '415.42 in Excel) so numeric-looking
    # prices aren't silently reinterpreted as floating point numbers.
    $cell.Value = "'" + $text
}

foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = $eValues[$row]
}
